# Generate Report for Archive
# - Update the localization status text from "Ready for handoff" to
#   "In Translation" everywhere it appears (Overview!E2, Overview!F2,
#   zh-cn!C2, de-de!C2 all share the same status string).
# - Re-fit the width of the columns that hold that status text on each
#   sheet now that the text is shorter.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview sheet keeps one status cell per locale column (E = zh-cn, F = de-de)
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# Per-locale detail sheets keep the status in column C
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# Resize the status columns to fit the new, shorter text
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
